# Auto-generated edit script: updates FFXIV Hyperion_Profits market-price
# snapshot values (currentAveragePrice / Leve profit columns H-N) across
# sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR per the scheduled runner diff.
$wb = $excel.ActiveWorkbook

# --- ALC (Worksheets index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H64").Value = 8155.364
$ws.Range("N64").Value = -8985.625
$ws.Range("M64").Value = -7016
$ws.Range("L64").Value = 8489.625
$ws.Range("I64").Value = 7264
$ws.Range("K64").Value = 7264
$ws.Range("J64").Value = 8489.625
$ws.Range("J67").Value = 8489.625
$ws.Range("H67").Value = 8155.364
$ws.Range("N67").Value = -10205.625
$ws.Range("M67").Value = -6406
$ws.Range("I67").Value = 7264
$ws.Range("L67").Value = 8489.625
$ws.Range("K67").Value = 7264
$ws.Range("J86").Value = 1621.6
$ws.Range("H86").Value = 1570.2
$ws.Range("N86").Value = -3867.6
$ws.Range("M86").Value = -421.5
$ws.Range("L86").Value = 1621.6
$ws.Range("I86").Value = 1544.5
$ws.Range("K86").Value = 1544.5
$ws.Range("J88").Value = 2333.7144
$ws.Range("H88").Value = 1929.3684
$ws.Range("N88").Value = -3145.7144
$ws.Range("L88").Value = 2333.7144
$ws.Range("I89").Value = 1544.5
$ws.Range("K89").Value = 7722.5
$ws.Range("H89").Value = 1570.2
$ws.Range("J89").Value = 1621.6
$ws.Range("N89").Value = -19340
$ws.Range("M89").Value = -2106.5
$ws.Range("L89").Value = 8108
$ws.Range("L91").Value = 2333.7144
$ws.Range("J91").Value = 2333.7144
$ws.Range("H91").Value = 1929.3684
$ws.Range("N91").Value = -5141.7144
$ws.Range("K92").Value = 2398.2856
$ws.Range("J92").Value = 749
$ws.Range("H92").Value = 2031.7778
$ws.Range("N92").Value = -3245
$ws.Range("M92").Value = -1150.2856
$ws.Range("L92").Value = 749
$ws.Range("I92").Value = 2398.2856
$ws.Range("J99").Value = 500
$ws.Range("H99").Value = 439.66666
$ws.Range("N99").Value = -4496
$ws.Range("M99").Value = 269.5
$ws.Range("L99").Value = 1500
$ws.Range("I99").Value = 409.5
$ws.Range("K99").Value = 1228.5
$ws.Range("M106").Value = -3368
$ws.Range("I106").Value = 3999
$ws.Range("K106").Value = 3999
$ws.Range("H106").Value = 3999

# --- ARM (Worksheets index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("J32").Value = 12179.777
$ws.Range("H32").Value = 3398.2307
$ws.Range("N32").Value = -12753.777
$ws.Range("M32").Value = -1699.9108
$ws.Range("I32").Value = 1986.9108
$ws.Range("L32").Value = 12179.777
$ws.Range("K32").Value = 1986.9108
$ws.Range("K102").Value = 4168246.8
$ws.Range("H102").Value = 3089069
$ws.Range("M102").Value = -4166624.8
$ws.Range("I102").Value = 4168246.8
$ws.Range("H122").Value = 870881.2
$ws.Range("M122").Value = -5717.749899999999
$ws.Range("I122").Value = 2722.5833
$ws.Range("K122").Value = 8167.749899999999

# --- BSM (Worksheets index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("I134").Value = 1268.2222
$ws.Range("K134").Value = 3804.6666
$ws.Range("H134").Value = 2825.024
$ws.Range("J134").Value = 12165.833
$ws.Range("N134").Value = -41567.499
$ws.Range("M134").Value = -1269.6666
$ws.Range("L134").Value = 36497.499
$ws.Range("L135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("H135").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N135").Value = $null
$ws.Range("N140").Value = $null

# --- CRP (Worksheets index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("M58").Value = -1911.4
$ws.Range("I58").Value = 2114.4
$ws.Range("K58").Value = 2114.4
$ws.Range("H58").Value = 2529.3794
$ws.Range("I62").Value = 2712.5715
$ws.Range("K62").Value = 2712.5715
$ws.Range("H62").Value = 2776.4443
$ws.Range("M62").Value = -2088.5715
$ws.Range("K65").Value = 13562.8575
$ws.Range("I65").Value = 2712.5715
$ws.Range("H65").Value = 2776.4443
$ws.Range("M65").Value = -10442.8575
$ws.Range("J99").Value = 4799
$ws.Range("H99").Value = 4749.375
$ws.Range("N99").Value = -7795
$ws.Range("M99").Value = -3168.6665
$ws.Range("L99").Value = 4799
$ws.Range("I99").Value = 4666.6665
$ws.Range("K99").Value = 4666.6665
$ws.Range("L126").Value = 14397
$ws.Range("I126").Value = 4666.6665
$ws.Range("K126").Value = 13999.9995
$ws.Range("J126").Value = 4799
$ws.Range("H126").Value = 4749.375
$ws.Range("N126").Value = -19337
$ws.Range("M126").Value = -11529.9995
$ws.Range("M132").Value = -180025.82
$ws.Range("L132").Value = 2666664
$ws.Range("I132").Value = 60851.94
$ws.Range("K132").Value = 182555.82
$ws.Range("J132").Value = 888888
$ws.Range("H132").Value = 106853.945
$ws.Range("N132").Value = -2671724
$ws.Range("I134").Value = 1581.2354
$ws.Range("K134").Value = 4743.706200000001
$ws.Range("H134").Value = 2219.05
$ws.Range("J134").Value = 5833.3335
$ws.Range("N134").Value = -22570.0005
$ws.Range("M134").Value = -2208.706200000001
$ws.Range("L134").Value = 17500.0005
$ws.Range("L135").Value = 115035.1
$ws.Range("J135").Value = 115035.1
$ws.Range("H135").Value = 115035.1
$ws.Range("N135").Value = -125175.1
$ws.Range("I136").Value = 2114.4
$ws.Range("K136").Value = 6343.200000000001
$ws.Range("H136").Value = 2529.3794
$ws.Range("M136").Value = -3793.200000000001
$ws.Range("J140").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

# --- CUL (Worksheets index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("I36").Value = 2
$ws.Range("K36").Value = 6
$ws.Range("H36").Value = 2
$ws.Range("M36").Value = 163
$ws.Range("J86").Value = 400
$ws.Range("H86").Value = 400
$ws.Range("N86").Value = -3572
$ws.Range("L86").Value = 1200
$ws.Range("H89").Value = 400
$ws.Range("J89").Value = 400
$ws.Range("N89").Value = -15456
$ws.Range("L89").Value = 3600
$ws.Range("M132").Value = -1293454.52
$ws.Range("L132").Value = 26496.999
$ws.Range("I132").Value = 143998.28
$ws.Range("K132").Value = 1295984.52
$ws.Range("J132").Value = 2944.111
$ws.Range("H132").Value = 64655.312
$ws.Range("N132").Value = -31556.999

# --- GSM (Worksheets index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("K2").Value = 1250.3334
$ws.Range("J2").Value = 8.75
$ws.Range("H2").Value = 868.3077
$ws.Range("N2").Value = -234.75
$ws.Range("M2").Value = -1137.3334
$ws.Range("L2").Value = 8.75
$ws.Range("I2").Value = 1250.3334
$ws.Range("J122").Value = 4936
$ws.Range("H122").Value = 638818.5600000001
$ws.Range("N122").Value = -19708
$ws.Range("M122").Value = -2674664.8
$ws.Range("L122").Value = 14808
$ws.Range("I122").Value = 892371.6
$ws.Range("K122").Value = 2677114.8
$ws.Range("I126").Value = 1821287.8
$ws.Range("K126").Value = 5463863.4
$ws.Range("H126").Value = 3791068.2
$ws.Range("M126").Value = -5461393.4

# --- LTW (Worksheets index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("L7").Value = 6614
$ws.Range("J7").Value = 6614
$ws.Range("H7").Value = 3887.375
$ws.Range("N7").Value = -6838
$ws.Range("I40").Value = 4812.4287
$ws.Range("K40").Value = 4812.4287
$ws.Range("H40").Value = 6283.3076
$ws.Range("M40").Value = -4676.4287
$ws.Range("K82").Value = 6173934
$ws.Range("H82").Value = 3969729.2
$ws.Range("J82").Value = 2161.2
$ws.Range("N82").Value = -2883.2
$ws.Range("M82").Value = -6173573
$ws.Range("L82").Value = 2161.2
$ws.Range("I82").Value = 6173934
$ws.Range("I85").Value = 6173934
$ws.Range("K85").Value = 6173934
$ws.Range("J85").Value = 2161.2
$ws.Range("H85").Value = 3969729.2
$ws.Range("N85").Value = -4657.2
$ws.Range("L85").Value = 2161.2
$ws.Range("M85").Value = -6172686
$ws.Range("H122").Value = 7107.6
$ws.Range("M122").Value = -9023.5
$ws.Range("I122").Value = 3824.5
$ws.Range("K122").Value = 11473.5
$ws.Range("L126").Value = 19842
$ws.Range("J126").Value = 6614
$ws.Range("H126").Value = 3887.375
$ws.Range("N126").Value = -24782
$ws.Range("J138").Value = 88000
$ws.Range("H138").Value = 88000
$ws.Range("N138").Value = -98280
$ws.Range("L138").Value = 88000

# --- WVR (Worksheets index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("L107").Value = 750
$ws.Range("I107").Value = 55556780
$ws.Range("K107").Value = 166670340
$ws.Range("J107").Value = 250
$ws.Range("H107").Value = 52632750
$ws.Range("N107").Value = -4590
$ws.Range("M107").Value = -166668420
$ws.Range("H122").Value = 1789.85
$ws.Range("M122").Value = -1750.3531
$ws.Range("I122").Value = 1400.1177
$ws.Range("K122").Value = 4200.3531
$ws.Range("M132").Value = -200024162
$ws.Range("I132").Value = 66675564
$ws.Range("K132").Value = 200026692
$ws.Range("H132").Value = 56168630
